$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 data: function name, filename, and line number
$ws.Range("A4").Value = "persistResolution"
$ws.Range("C4").Value = "/home/rdkv-core/cov/cov-analysis-linux64-2023.6.0/bin/device/devicesettings/persistResolution/generic/rpc/srv/dsVideoPort.c"
$ws.Range("D4").Value = 1636

# Update the active cell selection to D9
$ws.Range("D9").Select()
